$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (strikeout) values replacing the old Strike# values in column G
$kValues = @{
    2  = 3
    3  = 0
    4  = 3
    5  = 1
    6  = 3
    7  = 2
    8  = 3
    9  = 2
    10 = 3
    11 = 5
    12 = 1
    13 = 6
    14 = 4
    15 = 2
    16 = 6
    17 = 3
    18 = 8
    19 = 7
    20 = 5
    21 = 5
    22 = 3
    23 = 0
    24 = 5
    25 = 1
    26 = 7
    27 = 3
    28 = 6
    29 = 4
    30 = 1
    31 = 7
    32 = 4
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
